# Apply updates for trade #24 closed at 2026-02-17 13:18:37

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.7
$summary.Range("B6").Value = 24
$summary.Range("B9").Value = 29.17

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 24
$status.Range("G4").Value = 29.17

# --- New trade row data (used for both "All Trades" and "MarketMaking" sheets) ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(25, 1).Value = 24
    $ws.Cells.Item(25, 2).NumberFormat = "@"
    $ws.Cells.Item(25, 2).Value = "2026-02-17"
    $ws.Cells.Item(25, 3).NumberFormat = "@"
    $ws.Cells.Item(25, 3).Value = "13:18:37"
    $ws.Cells.Item(25, 4).Value = "MarketMaking"
    $ws.Cells.Item(25, 5).Value = "DOWN"
    $ws.Cells.Item(25, 6).Value = 0.8
    $ws.Cells.Item(25, 7).Value = 0.8
    $ws.Cells.Item(25, 8).Value = "CLOSED"
    $ws.Cells.Item(25, 9).Value = 0
    $ws.Cells.Item(25, 10).Value = 0
    $ws.Cells.Item(25, 11).Value = 99.17
    $ws.Cells.Item(25, 12).Value = 0
    $ws.Cells.Item(25, 13).Value = 0
    $ws.Cells.Item(25, 14).Value = 0.6
    $ws.Cells.Item(25, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(25, 16).Value = "early_exit"
    $ws.Cells.Item(25, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
